$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.626.83'
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").Value = '1.913.31'
$ws.Range("E3").Value = '  +3.98%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '''245.45'
$ws.Range("E5").Value = '  +5.64%  '
$ws.Range("D6").Value = '''0.633'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").Value = '''42.45'
$ws.Range("E8").Value = '  +3.77%  '
$ws.Range("D9").Value = '''0.339'
$ws.Range("E9").Value = '  +3.70%  '
$ws.Range("D10").Value = '''0.0709'
$ws.Range("E10").Value = '  +2.67%  '
$ws.Range("D11").Value = '''0.0997'
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").Value = '2.185.50'
$ws.Range("E12").Value = '  +3.73%  '
$ws.Range("E13").Value = '  +10.39%  '
$ws.Range("D14").Value = '1.919.33'
$ws.Range("E14").Value = '  +4.16%  '
$ws.Range("D15").Value = '''0.696'
$ws.Range("E15").Value = '  +3.91%  '
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("D17").Value = '35.611.02'
$ws.Range("E17").Value = '  +1.95%  '
$ws.Range("D18").Value = '''72.17'
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").Value = '''244.30'
$ws.Range("E20").Value = '  +1.85%  '
$ws.Range("E21").Value = '  +3.43%  '
$ws.Range("D22").Value = '''4.93'
$ws.Range("E22").Value = '  +4.11%  '
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").Value = '''2.28'
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("D25").Value = '''171.45'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").Value = '''2.12'
$ws.Range("E26").Value = '  +28.10%  '
$ws.Range("D27").Value = '''8.50'
$ws.Range("E27").Value = '  +8.27%  '
$ws.Range("D28").Value = '''18.02'
$ws.Range("E28").Value = '  +3.70%  '
$ws.Range("E29").Value = '  +2.09%  '
$ws.Range("D30").Value = '''4.11'
$ws.Range("E30").Value = '  +4.23%  '
$ws.Range("E31").Value = '  +2.94%  '
$ws.Range("D32").Value = '''0.946'
$ws.Range("E32").Value = '  +27.00%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").Value = '''4.15'
$ws.Range("D35").Value = '''1.74'
$ws.Range("E35").Value = '  +7.17%  '
$ws.Range("D36").Value = '''2.05'
$ws.Range("E36").Value = '  +5.52%  '
$ws.Range("E37").Value = '  +5.71%  '
$ws.Range("E38").Value = '  +5.22%  '
$ws.Range("E39").Value = '  +4.65%  '
$ws.Range("D40").Value = '''91.85'
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("D41").Value = '1.362.67'
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("D42").Value = '''15.25'
$ws.Range("E42").Value = '  +4.86%  '
$ws.Range("D43").Value = '''0.0595'
$ws.Range("E43").Value = '  +12.05%  '
$ws.Range("D44").Value = '''48.59'
$ws.Range("E44").Value = '  +43.06%  '
$ws.Range("B45").Value = 'Gas'
$ws.Range("C45").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D45").Value = '''13.08'
$ws.Range("E45").Value = '  +22.12%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''2.36'
$ws.Range("E46").Value = '  +4.96%  '
$ws.Range("E47").Value = '  +1.00%  '
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("E49").Value = '  +5.97%  '
$ws.Range("D50").Value = '2.094.18'
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("D51").Value = '''3.55'
$ws.Range("E51").Value = '  +5.00%  '
